# Commit: "Commit from C Drive"
#
# Applies to GitRepositorySetupNotes.docx:
#   1. Insert a new leading paragraph "Imp Doc from C Drive" and move the
#      _GoBack bookmark there (it is removed from its old home, the very
#      last paragraph of the document).
#   2. Mark "So" / "new" as grammar-flagged (w:proofErr gramStart/gramEnd)
#      by splitting their host runs.
#   3. Mark "AutomationTestNG" / "init" as spell-flagged
#      (w:proofErr spellStart/spellEnd) by splitting their host runs.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1a. Remove the _GoBack bookmark from its current location (last
#     paragraph of the body, just before the sectPr).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 1b. Insert the new first paragraph, with the run of text plus the
#     (re-created) _GoBack bookmark sitting right after it.
# ------------------------------------------------------------------
$introXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Imp Doc from C Drive</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$startRange = $d.Range(0, 0)
$null = $startRange.InsertXML($introXml)

# ------------------------------------------------------------------
# 2. "So my repository" -> "So" flagged gramStart/gramEnd.
# ------------------------------------------------------------------
$soXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005B63A7" w:rsidRDefault="005B63A7"><w:proofErr w:type="gramStart"/><w:r><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> my repository</w:t></w:r></w:p>
'@
$rng = $d.Content
$null = $rng.Find.Execute("So my repository", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pRange = $rng.Paragraphs(1).Range
$null = $pRange.InsertXML($soXml)

# ------------------------------------------------------------------
# 3. "To create A new Repository from command line" -> "new" flagged
#    gramStart/gramEnd.
# ------------------------------------------------------------------
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005B63A7" w:rsidRDefault="005B63A7"><w:r><w:t xml:space="preserve">To create A </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>new</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Repository from command line</w:t></w:r></w:p>
'@
$rng = $d.Content
$null = $rng.Find.Execute("To create A new Repository from command line", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pRange = $rng.Paragraphs(1).Range
$null = $pRange.InsertXML($newXml)

# ------------------------------------------------------------------
# 4. echo "# AutomationTestNG" >> README.md -> "AutomationTestNG"
#    flagged spellStart/spellEnd.
# ------------------------------------------------------------------
$echoXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005B63A7" w:rsidRPr="005B63A7" w:rsidRDefault="005B63A7" w:rsidP="005B63A7"><w:pPr><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005B63A7"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">echo "# </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005B63A7"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>AutomationTestNG</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005B63A7"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>" &gt;&gt; README.md</w:t></w:r></w:p>
'@
$rng = $d.Content
$null = $rng.Find.Execute('echo "# AutomationTestNG" >> README.md', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pRange = $rng.Paragraphs(1).Range
$null = $pRange.InsertXML($echoXml)

# ------------------------------------------------------------------
# 5. git init -> "init" flagged spellStart/spellEnd.
# ------------------------------------------------------------------
$gitInitXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005B63A7" w:rsidRPr="005B63A7" w:rsidRDefault="005B63A7" w:rsidP="005B63A7"><w:pPr><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="005B63A7"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">git </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005B63A7"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>init</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$rng = $d.Content
$null = $rng.Find.Execute("git init", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pRange = $rng.Paragraphs(1).Range
$null = $pRange.InsertXML($gitInitXml)

Write-Output "edits applied"
